$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4141732146861696
$ws.Range("C2").Value = 1.605476919376253
$ws.Range("D2").Value = 4.285373823189503
$ws.Range("E2").Value = 2.070114446881984
$ws.Range("F2").Value = 2.050671332274695
$ws.Range("G2").Value = 46
$ws.Range("B3").Value = -0.01716789264613588
$ws.Range("C3").Value = 1.651420560915961
$ws.Range("D3").Value = 4.541767649565558
$ws.Range("E3").Value = 2.131142334421978
$ws.Range("F3").Value = 2.155153870726067
$ws.Range("G3").Value = 45
$ws.Range("B4").Value = 0.4406255020385922
$ws.Range("C4").Value = 1.625614656034232
$ws.Range("D4").Value = 4.490980664127665
$ws.Range("E4").Value = 2.119193399415841
$ws.Range("F4").Value = 2.096844320931413
$ws.Range("G4").Value = 44
$ws.Range("B5").Value = 0.03988787390244446
$ws.Range("C5").Value = 1.807057129213736
$ws.Range("D5").Value = 5.145130224049834
$ws.Range("E5").Value = 2.268287949985591
$ws.Range("F5").Value = 2.294777636324286
$ws.Range("G5").Value = 43
$ws.Range("B6").Value = 0.455682680408077
$ws.Range("C6").Value = 1.62477962444376
$ws.Range("D6").Value = 4.358939589200301
$ws.Range("E6").Value = 2.087807364006627
$ws.Range("F6").Value = 2.062169714142621
$ws.Range("G6").Value = 42
$ws.Range("B7").Value = 0.1777087621310496
$ws.Range("C7").Value = 1.698096300376475
$ws.Range("D7").Value = 4.364176648270135
$ws.Range("E7").Value = 2.08906118825422
$ws.Range("F7").Value = 2.107346945862329
$ws.Range("G7").Value = 41
$ws.Range("B8").Value = 0.5039136335280825
$ws.Range("C8").Value = 1.664940487861543
$ws.Range("D8").Value = 4.498869753835015
$ws.Range("E8").Value = 2.121053925253909
$ws.Range("F8").Value = 2.086572653845029
$ws.Range("G8").Value = 40
$ws.Range("B9").Value = 0.1430518963212211
$ws.Range("C9").Value = 1.71910815220499
$ws.Range("D9").Value = 4.463347462520725
$ws.Range("E9").Value = 2.112663594262164
$ws.Range("F9").Value = 2.135369196980537
$ws.Range("G9").Value = 39
$ws.Range("B10").Value = 0.5887292808226858
$ws.Range("C10").Value = 1.678991187388869
$ws.Range("D10").Value = 4.640656357710296
$ws.Range("E10").Value = 2.154218270675071
$ws.Range("F10").Value = 2.100026121338619
$ws.Range("G10").Value = 38
$ws.Range("B11").Value = 0.1084470123171574
$ws.Range("C11").Value = 1.756824810856007
$ws.Range("D11").Value = 4.640724073720614
$ws.Range("E11").Value = 2.15423398769043
$ws.Range("F11").Value = 2.181179871918736
$ws.Range("G11").Value = 37
